# Re-run of the AMC BUR testing workbook against the frozen 2024 file.
# 1) Fix the transposed shared-string label used by the "cvd" column.
# 2) Push through the refreshed Voluntary Turnover / Internal Fill Rate percentages
#    (Jan-FY, PY Actual/AOP/Commit-Forecast rows) on every division tab.
$wb = $excel.ActiveWorkbook

# --- Shared string fix: "Voluntary Turnover Professional" -> "Professional Voluntary Turnover" ---
# The label lives in column C ("cvd") for rows 2-4 on every division sheet.
foreach ($i in 1..11) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($row in 2..4) {
        $cell = $ws.Cells.Item($row, 3)
        if ($cell.Value2 -eq "Voluntary Turnover Professional") {
            $cell.Value = "Professional Voluntary Turnover"
        }
    }
}

# --- Sheet 1: AMC Aerospace Solutions Division ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 0.0229
$ws.Range("I2").Value = 0.0326
$ws.Range("J2").Value = 0.0181
$ws.Range("M2").Value = 0.0357
$ws.Range("Q2").Value = 0.0298
$ws.Range("R2").Value = 0.0162
$ws.Range("U2").Value = 0.0198
$ws.Range("V2").Value = 0.1165
$ws.Range("H3").Value = 0.02061
$ws.Range("I3").Value = 0.02934
$ws.Range("J3").Value = 0.01629
$ws.Range("M3").Value = 0.03213
$ws.Range("Q3").Value = 0.02682
$ws.Range("R3").Value = 0.01458
$ws.Range("U3").Value = 0.01782
$ws.Range("V3").Value = 0.10485
$ws.Range("M4").Value = 0.03735
$ws.Range("Q4").Value = 0.03735
$ws.Range("U4").Value = 0.03735

# --- Sheet 3: AMC Conveyance Solutions Division ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 0.0097
$ws.Range("G2").Value = 0.0077
$ws.Range("H2").Value = 0.0058
$ws.Range("I2").Value = 0.0232
$ws.Range("J2").Value = 0.0154
$ws.Range("K2").Value = 0.0136
$ws.Range("L2").Value = 0.0059
$ws.Range("M2").Value = 0.035
$ws.Range("N2").Value = 0.0061
$ws.Range("O2").Value = 0.0166
$ws.Range("P2").Value = 0.0209
$ws.Range("Q2").Value = 0.0432
$ws.Range("R2").Value = 0.0083
$ws.Range("S2").Value = 0.0063
$ws.Range("T2").Value = 0.0041
$ws.Range("U2").Value = 0.0187
$ws.Range("V2").Value = 0.1202
$ws.Range("F3").Value = 0.00873
$ws.Range("G3").Value = 0.00693
$ws.Range("H3").Value = 0.00522
$ws.Range("I3").Value = 0.02088
$ws.Range("J3").Value = 0.01386
$ws.Range("K3").Value = 0.01224
$ws.Range("L3").Value = 0.00531
$ws.Range("M3").Value = 0.0315
$ws.Range("N3").Value = 0.00549
$ws.Range("O3").Value = 0.01494
$ws.Range("P3").Value = 0.01881
$ws.Range("Q3").Value = 0.03888
$ws.Range("R3").Value = 0.00747
$ws.Range("S3").Value = 0.00567
$ws.Range("T3").Value = 0.00369
$ws.Range("U3").Value = 0.01683
$ws.Range("V3").Value = 0.10818
$ws.Range("M4").Value = 0.0333
$ws.Range("Q4").Value = 0.0333
$ws.Range("U4").Value = 0.0333

# --- Sheet 4: AMC India/Middle East ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("V2").Value = 0.1587
$ws.Range("V3").Value = 0.14283

# --- Sheet 5: AMC Linear Motion Division ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("F2").Value = 0.0101
$ws.Range("G2").Value = 0.004
$ws.Range("H2").Value = 0.0061
$ws.Range("I2").Value = 0.0203
$ws.Range("J2").Value = 0.0083
$ws.Range("K2").Value = 0.0146
$ws.Range("L2").Value = 0.0063
$ws.Range("M2").Value = 0.0292
$ws.Range("N2").Value = 0.0084
$ws.Range("O2").Value = 0.0132
$ws.Range("P2").Value = 0.0066
$ws.Range("Q2").Value = 0.0281
$ws.Range("R2").Value = 0.0111
$ws.Range("S2").Value = 0.0067
$ws.Range("T2").Value = 0.0067
$ws.Range("U2").Value = 0.0246
$ws.Range("V2").Value = 0.102
$ws.Range("F3").Value = 0.00909
$ws.Range("G3").Value = 0.0036
$ws.Range("H3").Value = 0.00549
$ws.Range("I3").Value = 0.01827
$ws.Range("J3").Value = 0.00747
$ws.Range("K3").Value = 0.01314
$ws.Range("L3").Value = 0.00567
$ws.Range("M3").Value = 0.02628
$ws.Range("N3").Value = 0.00756
$ws.Range("O3").Value = 0.01188
$ws.Range("P3").Value = 0.00594
$ws.Range("Q3").Value = 0.02529
$ws.Range("R3").Value = 0.00999
$ws.Range("S3").Value = 0.00603
$ws.Range("T3").Value = 0.00603
$ws.Range("U3").Value = 0.02214
$ws.Range("V3").Value = 0.0918
$ws.Range("M4").Value = 0.023025
$ws.Range("Q4").Value = 0.023025
$ws.Range("U4").Value = 0.023025

# --- Sheet 6: AMC Micro-Motion Division ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("F2").Value = 0.0086
$ws.Range("G2").Value = 0.0203
$ws.Range("H2").Value = 0.0087
$ws.Range("I2").Value = 0.0376
$ws.Range("J2").Value = 0.0144
$ws.Range("K2").Value = 0.0203
$ws.Range("L2").Value = 0.0116
$ws.Range("M2").Value = 0.0463
$ws.Range("Q2").Value = 0.0197
$ws.Range("R2").Value = 0.0111
$ws.Range("S2").Value = 0.0084
$ws.Range("T2").Value = 0.0139
$ws.Range("U2").Value = 0.0333
$ws.Range("V2").Value = 0.1364
$ws.Range("F3").Value = 0.00774
$ws.Range("G3").Value = 0.01827
$ws.Range("H3").Value = 0.00783
$ws.Range("I3").Value = 0.03384
$ws.Range("J3").Value = 0.01296
$ws.Range("K3").Value = 0.01827
$ws.Range("L3").Value = 0.01044
$ws.Range("M3").Value = 0.04167
$ws.Range("Q3").Value = 0.01773
$ws.Range("R3").Value = 0.00999
$ws.Range("S3").Value = 0.00756
$ws.Range("T3").Value = 0.01251
$ws.Range("U3").Value = 0.02997
$ws.Range("V3").Value = 0.12276
$ws.Range("M4").Value = 0.0189
$ws.Range("Q4").Value = 0.0189
$ws.Range("U4").Value = 0.0189

# --- Sheet 7: AMC Motion Control Systems Division ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("U2").Value = 0.0274
$ws.Range("U3").Value = 0.02466
$ws.Range("M4").Value = 0.0282
$ws.Range("Q4").Value = 0.0282
$ws.Range("U4").Value = 0.0282

# --- Sheet 8: AMC Power Management Division ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("F2").Value = 0.008
$ws.Range("G2").Value = 0.0078
$ws.Range("H2").Value = 0.0078
$ws.Range("I2").Value = 0.0236
$ws.Range("K2").Value = 0.015
$ws.Range("L2").Value = 0.0148
$ws.Range("M2").Value = 0.03
$ws.Range("N2").Value = 0.0074
$ws.Range("O2").Value = 0.0074
$ws.Range("P2").Value = 0.0226
$ws.Range("Q2").Value = 0.0372
$ws.Range("S2").Value = 0.0147
$ws.Range("T2").Value = 0.0147
$ws.Range("U2").Value = 0.0295
$ws.Range("V2").Value = 0.1206
$ws.Range("F3").Value = 0.0072
$ws.Range("G3").Value = 0.00702
$ws.Range("H3").Value = 0.00702
$ws.Range("I3").Value = 0.02124
$ws.Range("K3").Value = 0.0135
$ws.Range("L3").Value = 0.01332
$ws.Range("M3").Value = 0.027
$ws.Range("N3").Value = 0.00666
$ws.Range("O3").Value = 0.00666
$ws.Range("P3").Value = 0.02034
$ws.Range("Q3").Value = 0.03348
$ws.Range("S3").Value = 0.01323
$ws.Range("T3").Value = 0.01323
$ws.Range("U3").Value = 0.02655
$ws.Range("V3").Value = 0.10854
$ws.Range("M4").Value = 0.0219
$ws.Range("Q4").Value = 0.0219
$ws.Range("U4").Value = 0.0219

# --- Sheet 9: AMC Segment Functions ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("M4").Value = 0.04545
$ws.Range("Q4").Value = 0.04545
$ws.Range("U4").Value = 0.04545

# --- Sheet 10: AMC Thomson Linear Motion - General ---
$ws = $wb.Worksheets.Item(10)
$ws.Range("N2").Value = 0.0139
$ws.Range("O2").Value = 0.0299
$ws.Range("Q2").Value = 0.0435
$ws.Range("R2").Value = 0.0147
$ws.Range("S2").Value = 0.0152
$ws.Range("T2").Value = 0.0308
$ws.Range("U2").Value = 0.0603
$ws.Range("V2").Value = 0.1029
$ws.Range("N3").Value = 0.01251
$ws.Range("O3").Value = 0.02691
$ws.Range("Q3").Value = 0.03915
$ws.Range("R3").Value = 0.01323
$ws.Range("S3").Value = 0.01368
$ws.Range("T3").Value = 0.02772
$ws.Range("U3").Value = 0.05427
$ws.Range("V3").Value = 0.09261
$ws.Range("M4").Value = 0.75
$ws.Range("Q4").Value = 0.75
$ws.Range("U4").Value = 0.75

# --- Sheet 11: L1_AMC ---
$ws = $wb.Worksheets.Item(11)
$ws.Range("F2").Value = 0.0076
$ws.Range("G2").Value = 0.0057
$ws.Range("H2").Value = 0.0091
$ws.Range("I2").Value = 0.0224
$ws.Range("J2").Value = 0.0092
$ws.Range("K2").Value = 0.0107
$ws.Range("L2").Value = 0.008
$ws.Range("M2").Value = 0.0279
$ws.Range("N2").Value = 0.007
$ws.Range("O2").Value = 0.0114
$ws.Range("P2").Value = 0.0106
$ws.Range("Q2").Value = 0.0289
$ws.Range("R2").Value = 0.009
$ws.Range("S2").Value = 0.0075
$ws.Range("U2").Value = 0.026
$ws.Range("V2").Value = 0.1051
$ws.Range("F3").Value = 0.00684
$ws.Range("G3").Value = 0.00513
$ws.Range("H3").Value = 0.00819
$ws.Range("I3").Value = 0.02016
$ws.Range("J3").Value = 0.00828
$ws.Range("K3").Value = 0.00963
$ws.Range("L3").Value = 0.0072
$ws.Range("M3").Value = 0.02511
$ws.Range("N3").Value = 0.0063
$ws.Range("O3").Value = 0.01026
$ws.Range("P3").Value = 0.00954
$ws.Range("Q3").Value = 0.02601
$ws.Range("R3").Value = 0.0081
$ws.Range("S3").Value = 0.00675
$ws.Range("U3").Value = 0.0234
$ws.Range("V3").Value = 0.09459
$ws.Range("M4").Value = 0.02745
$ws.Range("Q4").Value = 0.02745
$ws.Range("U4").Value = 0.02745
